$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Delete()

$ws.Range("B2").Value = "https://casa.mercadolibre.com.mx/MLM-2021698087-residencia-de-autor-con-acabados-premium-altozano-_JM#position=1&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://casa.mercadolibre.com.mx/MLM-2021698087-residencia-de-autor-con-acabados-premium-altozano-_JM", "position=1&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B2").Style = "Hyperlink"

$ws.Range("B3").Value = "https://casa.mercadolibre.com.mx/MLM-2021787853-vivir-con-altura-altozano-_JM#position=2&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B3"), "https://casa.mercadolibre.com.mx/MLM-2021787853-vivir-con-altura-altozano-_JM", "position=2&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B3").Style = "Hyperlink"

$ws.Range("B4").Value = "https://casa.mercadolibre.com.mx/MLM-2020248527-casa-condominio-en-venta-santa-ines-lomas-el-campanario-_JM#position=3&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://casa.mercadolibre.com.mx/MLM-2020248527-casa-condominio-en-venta-santa-ines-lomas-el-campanario-_JM", "position=3&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B4").Style = "Hyperlink"

$ws.Range("B5").Value = "https://casa.mercadolibre.com.mx/MLM-2023177877-casa-en-venta-en-torre-de-piedra-gran-reserva-con-3-habitac-_JM#position=4&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B5"), "https://casa.mercadolibre.com.mx/MLM-2023177877-casa-en-venta-en-torre-de-piedra-gran-reserva-con-3-habitac-_JM", "position=4&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B5").Style = "Hyperlink"

$ws.Range("B6").Value = "https://casa.mercadolibre.com.mx/MLM-2023037855-casa-en-venta-en-altozano-de-2-recamaras-moderna-en-priv-_JM#position=5&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B6"), "https://casa.mercadolibre.com.mx/MLM-2023037855-casa-en-venta-en-altozano-de-2-recamaras-moderna-en-priv-_JM", "position=5&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("B7").Value = "https://casa.mercadolibre.com.mx/MLM-2890670236-casa-en-venta-en-zibata-de-4-recamaras-y-roof-garden-ideal-_JM#position=6&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://casa.mercadolibre.com.mx/MLM-2890670236-casa-en-venta-en-zibata-de-4-recamaras-y-roof-garden-ideal-_JM", "position=6&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B7").Style = "Hyperlink"

$ws.Range("B8").Value = "https://casa.mercadolibre.com.mx/MLM-2023071019-casa-en-venta-en-zibata-con-4-habitaciones-y-cuarto-de-serv-_JM#position=7&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B8"), "https://casa.mercadolibre.com.mx/MLM-2023071019-casa-en-venta-en-zibata-con-4-habitaciones-y-cuarto-de-serv-_JM", "position=7&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B8").Style = "Hyperlink"

$ws.Range("B9").Value = "https://casa.mercadolibre.com.mx/MLM-2023011467-casa-en-venta-en-zibata-estilo-mexico-contemporanea-de-3-_JM#position=8&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://casa.mercadolibre.com.mx/MLM-2023011467-casa-en-venta-en-zibata-estilo-mexico-contemporanea-de-3-_JM", "position=8&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B9").Style = "Hyperlink"

$ws.Range("B10").Value = "https://casa.mercadolibre.com.mx/MLM-2889893070-4-habitaciones-amplio-jardin-la-condesa-juriquilla-_JM#position=9&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B10"), "https://casa.mercadolibre.com.mx/MLM-2889893070-4-habitaciones-amplio-jardin-la-condesa-juriquilla-_JM", "position=9&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B10").Style = "Hyperlink"

$ws.Range("B11").Value = "https://casa.mercadolibre.com.mx/MLM-2889869734-4-habitaciones-amplio-jardin-la-condesa-juriquilla-_JM#position=10&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B11"), "https://casa.mercadolibre.com.mx/MLM-2889869734-4-habitaciones-amplio-jardin-la-condesa-juriquilla-_JM", "position=10&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B11").Style = "Hyperlink"

$ws.Range("B12").Value = "https://casa.mercadolibre.com.mx/MLM-2021853869-residencia-con-acabados-premium-altozano-_JM#position=11&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B12"), "https://casa.mercadolibre.com.mx/MLM-2021853869-residencia-con-acabados-premium-altozano-_JM", "position=11&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B12").Style = "Hyperlink"

$ws.Range("B13").Value = "https://casa.mercadolibre.com.mx/MLM-2021698085-residencia-con-amplios-espacios-jurica-_JM#position=12&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B13"), "https://casa.mercadolibre.com.mx/MLM-2021698085-residencia-con-amplios-espacios-jurica-_JM", "position=12&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B13").Style = "Hyperlink"

$ws.Range("B14").Value = "https://casa.mercadolibre.com.mx/MLM-2021800157-increible-casa-de-autor-zibata-_JM#position=13&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B14"), "https://casa.mercadolibre.com.mx/MLM-2021800157-increible-casa-de-autor-zibata-_JM", "position=13&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B14").Style = "Hyperlink"

$ws.Range("B15").Value = "https://casa.mercadolibre.com.mx/MLM-2021800155-doble-altura-amplio-jardin-la-vista-residencial-_JM#position=14&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B15"), "https://casa.mercadolibre.com.mx/MLM-2021800155-doble-altura-amplio-jardin-la-vista-residencial-_JM", "position=14&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B15").Style = "Hyperlink"

$ws.Range("B16").Value = "https://casa.mercadolibre.com.mx/MLM-2021762251-amplio-jardin-con-terraza-la-vista-residencial-_JM#position=15&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B16"), "https://casa.mercadolibre.com.mx/MLM-2021762251-amplio-jardin-con-terraza-la-vista-residencial-_JM", "position=15&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B16").Style = "Hyperlink"

$ws.Range("B17").Value = "https://casa.mercadolibre.com.mx/MLM-2885801136-la-vida-que-deseas-altozano-_JM#position=16&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B17"), "https://casa.mercadolibre.com.mx/MLM-2885801136-la-vida-que-deseas-altozano-_JM", "position=16&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B17").Style = "Hyperlink"

$ws.Range("B18").Value = "https://casa.mercadolibre.com.mx/MLM-2021762245-exclusividad-y-confort-zibata-_JM#position=17&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B18"), "https://casa.mercadolibre.com.mx/MLM-2021762245-exclusividad-y-confort-zibata-_JM", "position=17&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B18").Style = "Hyperlink"

$ws.Range("B19").Value = "https://casa.mercadolibre.com.mx/MLM-2021749071-vive-en-plenitud-lomas-del-campanario-norte-_JM#position=18&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B19"), "https://casa.mercadolibre.com.mx/MLM-2021749071-vive-en-plenitud-lomas-del-campanario-norte-_JM", "position=18&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B19").Style = "Hyperlink"

$ws.Range("B20").Value = "https://casa.mercadolibre.com.mx/MLM-2885801130-nuevas-emociones-altozano-_JM#position=19&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B20"), "https://casa.mercadolibre.com.mx/MLM-2885801130-nuevas-emociones-altozano-_JM", "position=19&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B20").Style = "Hyperlink"

$ws.Range("B21").Value = "https://casa.mercadolibre.com.mx/MLM-2021710635-residencia-con-acabados-premium-altozano-_JM#position=20&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B21"), "https://casa.mercadolibre.com.mx/MLM-2021710635-residencia-con-acabados-premium-altozano-_JM", "position=20&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B21").Style = "Hyperlink"

$ws.Range("B22").Value = "https://casa.mercadolibre.com.mx/MLM-2021762239-una-vida-excepcional-altozano-_JM#position=21&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B22"), "https://casa.mercadolibre.com.mx/MLM-2021762239-una-vida-excepcional-altozano-_JM", "position=21&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B22").Style = "Hyperlink"

$ws.Range("B23").Value = "https://casa.mercadolibre.com.mx/MLM-2021774945-vive-en-armonia-canadas-del-arroyo-_JM#position=22&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B23"), "https://casa.mercadolibre.com.mx/MLM-2021774945-vive-en-armonia-canadas-del-arroyo-_JM", "position=22&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B23").Style = "Hyperlink"

$ws.Range("B24").Value = "https://casa.mercadolibre.com.mx/MLM-2021710631-la-vida-que-deseas-la-vista-residencial-_JM#position=23&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B24"), "https://casa.mercadolibre.com.mx/MLM-2021710631-la-vida-que-deseas-la-vista-residencial-_JM", "position=23&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B24").Style = "Hyperlink"

$ws.Range("B25").Value = "https://casa.mercadolibre.com.mx/MLM-2021774937-terraza-jardin-doble-altura-zibata-_JM#position=24&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B25"), "https://casa.mercadolibre.com.mx/MLM-2021774937-terraza-jardin-doble-altura-zibata-_JM", "position=24&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B25").Style = "Hyperlink"

$ws.Range("B26").Value = "https://casa.mercadolibre.com.mx/MLM-2021762233-amplio-jardin-habitacion-en-planta-baja-zibata-_JM#position=25&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B26"), "https://casa.mercadolibre.com.mx/MLM-2021762233-amplio-jardin-habitacion-en-planta-baja-zibata-_JM", "position=25&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B26").Style = "Hyperlink"

$ws.Range("B27").Value = "https://casa.mercadolibre.com.mx/MLM-2885801124-diseno-unico-y-moderno-zibata-_JM#position=26&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B27"), "https://casa.mercadolibre.com.mx/MLM-2885801124-diseno-unico-y-moderno-zibata-_JM", "position=26&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B27").Style = "Hyperlink"

$ws.Range("B28").Value = "https://casa.mercadolibre.com.mx/MLM-2885801122-terraza-jardin-roof-garden-zibata-_JM#position=27&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B28"), "https://casa.mercadolibre.com.mx/MLM-2885801122-terraza-jardin-roof-garden-zibata-_JM", "position=27&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B28").Style = "Hyperlink"

$ws.Range("B29").Value = "https://casa.mercadolibre.com.mx/MLM-2021762231-una-casa-de-altura-lomas-del-campanario-norte-_JM#position=28&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B29"), "https://casa.mercadolibre.com.mx/MLM-2021762231-una-casa-de-altura-lomas-del-campanario-norte-_JM", "position=28&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B29").Style = "Hyperlink"

$ws.Range("B30").Value = "https://casa.mercadolibre.com.mx/MLM-2885801116-conocela-y-enamorate-altozano-_JM#position=29&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B30"), "https://casa.mercadolibre.com.mx/MLM-2885801116-conocela-y-enamorate-altozano-_JM", "position=29&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B30").Style = "Hyperlink"

$ws.Range("B31").Value = "https://casa.mercadolibre.com.mx/MLM-2021736165-conocela-y-enamorate-lomas-del-campanario-norte-_JM#position=30&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B31"), "https://casa.mercadolibre.com.mx/MLM-2021736165-conocela-y-enamorate-lomas-del-campanario-norte-_JM", "position=30&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B31").Style = "Hyperlink"

$ws.Range("B32").Value = "https://casa.mercadolibre.com.mx/MLM-2021800131-jardin-amplia-estancia-acceso-a-roof-colinas-de-juriquill-_JM#position=31&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B32"), "https://casa.mercadolibre.com.mx/MLM-2021800131-jardin-amplia-estancia-acceso-a-roof-colinas-de-juriquill-_JM", "position=31&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B32").Style = "Hyperlink"

$ws.Range("B33").Value = "https://casa.mercadolibre.com.mx/MLM-2021736163-diseno-que-enamora-colinas-de-juriquilla-_JM#position=32&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B33"), "https://casa.mercadolibre.com.mx/MLM-2021736163-diseno-que-enamora-colinas-de-juriquilla-_JM", "position=32&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B33").Style = "Hyperlink"

$ws.Range("B34").Value = "https://casa.mercadolibre.com.mx/MLM-2021723537-canadas-del-arroyo-habitacion-planta-baja-_JM#position=33&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B34"), "https://casa.mercadolibre.com.mx/MLM-2021723537-canadas-del-arroyo-habitacion-planta-baja-_JM", "position=33&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B34").Style = "Hyperlink"

$ws.Range("B35").Value = "https://casa.mercadolibre.com.mx/MLM-2890670194-departamento-en-venta-en-central-park-amueblado-de-lujo-_JM#position=34&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B35"), "https://casa.mercadolibre.com.mx/MLM-2890670194-departamento-en-venta-en-central-park-amueblado-de-lujo-_JM", "position=34&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B35").Style = "Hyperlink"

$ws.Range("B36").Value = "https://casa.mercadolibre.com.mx/MLM-2890740674-casa-en-venta-en-zibata-con-3-recamaras-dentro-de-condomin-_JM#position=35&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B36"), "https://casa.mercadolibre.com.mx/MLM-2890740674-casa-en-venta-en-zibata-con-3-recamaras-dentro-de-condomin-_JM", "position=35&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B36").Style = "Hyperlink"

$ws.Range("B37").Value = "https://casa.mercadolibre.com.mx/MLM-2021710639-hasta-5-recamaras-roof-garden-milenio-lll-_JM#position=36&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B37"), "https://casa.mercadolibre.com.mx/MLM-2021710639-hasta-5-recamaras-roof-garden-milenio-lll-_JM", "position=36&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B37").Style = "Hyperlink"

$ws.Range("B38").Value = "https://casa.mercadolibre.com.mx/MLM-2890669990-casa-en-venta-en-santa-fe-juriquilla-dentro-de-privada-_JM#position=37&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B38"), "https://casa.mercadolibre.com.mx/MLM-2890669990-casa-en-venta-en-santa-fe-juriquilla-dentro-de-privada-_JM", "position=37&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B38").Style = "Hyperlink"

$ws.Range("B39").Value = "https://casa.mercadolibre.com.mx/MLM-2021723553-distribuida-mayormente-en-una-planta-jurica-_JM#position=38&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B39"), "https://casa.mercadolibre.com.mx/MLM-2021723553-distribuida-mayormente-en-una-planta-jurica-_JM", "position=38&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B39").Style = "Hyperlink"

$ws.Range("B40").Value = "https://casa.mercadolibre.com.mx/MLM-2021710637-acceso-a-roof-con-vista-pedregal-de-schoenstatt-_JM#position=39&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B40"), "https://casa.mercadolibre.com.mx/MLM-2021710637-acceso-a-roof-con-vista-pedregal-de-schoenstatt-_JM", "position=39&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B40").Style = "Hyperlink"

$ws.Range("B41").Value = "https://casa.mercadolibre.com.mx/MLM-2885801128-casa-ensueno-juriquilla-_JM#position=40&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B41"), "https://casa.mercadolibre.com.mx/MLM-2885801128-casa-ensueno-juriquilla-_JM", "position=40&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B41").Style = "Hyperlink"

$ws.Range("B42").Value = "https://casa.mercadolibre.com.mx/MLM-2021736177-habitacion-planta-baja-lomas-de-juriquilla-_JM#position=41&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B42"), "https://casa.mercadolibre.com.mx/MLM-2021736177-habitacion-planta-baja-lomas-de-juriquilla-_JM", "position=41&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B42").Style = "Hyperlink"

$ws.Range("B43").Value = "https://casa.mercadolibre.com.mx/MLM-2021736173-el-lujo-que-mereces-altozano-_JM#position=42&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B43"), "https://casa.mercadolibre.com.mx/MLM-2021736173-el-lujo-que-mereces-altozano-_JM", "position=42&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B43").Style = "Hyperlink"

$ws.Range("B44").Value = "https://casa.mercadolibre.com.mx/MLM-2021800145-diseno-distintivo-zibata-_JM#position=43&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B44"), "https://casa.mercadolibre.com.mx/MLM-2021800145-diseno-distintivo-zibata-_JM", "position=43&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B44").Style = "Hyperlink"

$ws.Range("B45").Value = "https://casa.mercadolibre.com.mx/MLM-2021787865-tu-nuevo-hogar-residencial-el-refugio-_JM#position=44&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B45"), "https://casa.mercadolibre.com.mx/MLM-2021787865-tu-nuevo-hogar-residencial-el-refugio-_JM", "position=44&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B45").Style = "Hyperlink"

$ws.Range("B46").Value = "https://casa.mercadolibre.com.mx/MLM-2890670252-casa-en-venta-en-altozano-dentro-de-condominio-con-alta-pl-_JM#position=45&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B46"), "https://casa.mercadolibre.com.mx/MLM-2890670252-casa-en-venta-en-altozano-dentro-de-condominio-con-alta-pl-_JM", "position=45&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B46").Style = "Hyperlink"

$ws.Range("B47").Value = "https://casa.mercadolibre.com.mx/MLM-2890753436-exclusiva-casa-en-venta-en-altozano-dentro-de-condominio-c-_JM#position=46&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B47"), "https://casa.mercadolibre.com.mx/MLM-2890753436-exclusiva-casa-en-venta-en-altozano-dentro-de-condominio-c-_JM", "position=46&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B47").Style = "Hyperlink"

$ws.Range("B48").Value = "https://casa.mercadolibre.com.mx/MLM-2021698081-iluminacion-y-amplitud-lomas-del-campanario-norte-_JM#position=47&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B48"), "https://casa.mercadolibre.com.mx/MLM-2021698081-iluminacion-y-amplitud-lomas-del-campanario-norte-_JM", "position=47&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B48").Style = "Hyperlink"

$ws.Range("B49").Value = "https://casa.mercadolibre.com.mx/MLM-2885801114-diseno-distintivo-la-espiga-_JM#position=48&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc"
$ws.Hyperlinks.Add($ws.Range("B49"), "https://casa.mercadolibre.com.mx/MLM-2885801114-diseno-distintivo-la-espiga-_JM", "position=48&search_layout=grid&type=item&tracking_id=fa943906-1423-4600-bb54-a94262abe0fc")
$ws.Range("B49").Style = "Hyperlink"
